$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ------------------------------------------------------------------
# 1) Insert two new paragraphs right after "git push origin master"
#    and before "TypeScript  commands :" :
#      - a plain (unformatted) paragraph
#      - a bold paragraph with the text "git pull"
# ------------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("git push origin master", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $anchorPara = $find1.Paragraphs(1)
    $anchorIndex = $anchorPara.Index
    $anchorRange = $anchorPara.Range
    $anchorRange.Collapse(0)
    $anchorRange.InsertParagraphAfter()
    $anchorRange.InsertParagraphAfter()

    $newPara1 = $d.Paragraphs($anchorIndex + 1)
    $newPara2 = $d.Paragraphs($anchorIndex + 2)

    $xmlPlain = "<w:p $wNs><w:r><w:t>git  pull   to update the checkout folder from git hub</w:t></w:r></w:p>"
    $xmlBold  = "<w:p $wNs><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>git pull</w:t></w:r></w:p>"

    $newPara1.Range.InsertXML($xmlPlain)
    $newPara2.Range.InsertXML($xmlBold)
}

# ------------------------------------------------------------------
# 2) Add <w:lastRenderedPageBreak/> before the text run in the
#    paragraph "First install the node -&gt;npm -&gt; typescript-&gt;angular"
# ------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("First install the node", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para2 = $find2.Paragraphs(1)
    $xml2 = "<w:p $wNs><w:r><w:lastRenderedPageBreak/><w:t>First install the node -&gt;npm -&gt; typescript-&gt;angular</w:t></w:r></w:p>"
    $para2.Range.InsertXML($xml2)
}

# ------------------------------------------------------------------
# 3) Remove <w:lastRenderedPageBreak/> from the paragraph
#    "npm   install -g  @angular/cli   (after this check with ng command like    - ng  version)"
# ------------------------------------------------------------------
$find3 = $d.Content
$found3 = $find3.Find.Execute("@angular/cli   (after this", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $para3 = $find3.Paragraphs(1)
    $dash = [char]0x2013
    $xml3 = "<w:p $wNs><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>npm   install " + $dash + "g  @angular/cli</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>   (after this check with ng command like    - ng  version)</w:t></w:r></w:p>"
    $para3.Range.InsertXML($xml3)
}

Write-Output "edits applied"
